$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "0.4.0-snapshot-1"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-05-23T12:16:26+00:00"
$ws.Range("B10").Value = "ANS (https://esante.gouv.fr)"
